# Generate Report for Handback
# Update the "Correspond Handoff/Handback Datetime" and "Latest HO Xliff Generate Date"
# values across the Overview, zh-cn and de-de worksheets to reflect the newest
# handback report generation times.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the 862cc7d4... row
$wsOverview.Range("G4").Value = "2016-11-08 22:43:49"

# zh-cn sheet: "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# for the 862cc7d4... row
$wsZhCn.Range("H4").Value = "2016-11-08 22:43:35"
$wsZhCn.Range("K4").Value = "2016-11-08 22:44:29"

# de-de sheet: "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# for the 862cc7d4... row
$wsDeDe.Range("H4").Value = "2016-11-08 22:43:49"
$wsDeDe.Range("K4").Value = "2016-11-08 22:44:46"
